$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A8").PasteSpecial(-4122) | Out-Null
$ws.Range("A8").Value = 44490.551210254634
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B8").PasteSpecial(-4122) | Out-Null
$ws.Range("B8").Value = 'Tyler Cowie'
$ws.Range("C7").Copy() | Out-Null
$ws.Range("C8").PasteSpecial(-4122) | Out-Null
$ws.Range("C8").Value = 'Innovation'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D8").PasteSpecial(-4122) | Out-Null
$ws.Range("D8").Value = 'Ontario Electricity Demand Forecasting'
$ws.Range("E7").Copy() | Out-Null
$ws.Range("E8").PasteSpecial(-4122) | Out-Null
$ws.Range("E8").Value = 'Yes'
$ws.Range("F7").Copy() | Out-Null
$ws.Range("F8").PasteSpecial(-4122) | Out-Null
$ws.Range("F8").Value = 'Website'
$ws.Range("G7").Copy() | Out-Null
$ws.Range("G8").PasteSpecial(-4122) | Out-Null
$ws.Range("G8").Value = 'Maybe'
$ws.Range("H7").Copy() | Out-Null
$ws.Range("H8").PasteSpecial(-4122) | Out-Null
$ws.Range("H8").Value = 'The goal of our project is to create a model to predict when the peak demand hours will be on the Ontario power grid.This model will be provided to Queen''s Energy Management to allow them to save money by limiting their grid power usage during these time periods. Currently, we are looking at using a logistic regression to find the probability of a peak demand hour, however we are also considering the use of multi-layer perceptron for complete grid forecasting.'
$ws.Range("I7").Copy() | Out-Null
$ws.Range("I8").PasteSpecial(-4122) | Out-Null
$ws.Range("I8").Value = 'Currently, there are other Ontario grid demand forecasting tools on the market, which take the form of online dashboards featuring many analytics. However, these tools cost companies thousands of dollars every month and are unattainable for smaller companies. We could potentially produce our own online dashboard for Ontario''s grid at a fraction of the cost, including other analytics not included in the leading competitor''s products (long-term forecasts, humidity, aggregation of other forecasting tools).'
$ws.Range("J7").Copy() | Out-Null
$ws.Range("J8").PasteSpecial(-4122) | Out-Null
$ws.Range("J8").Value = 'N/A'

# Row 9
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A9").PasteSpecial(-4122) | Out-Null
$ws.Range("A9").Value = 44491.011418333335
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B9").PasteSpecial(-4122) | Out-Null
$ws.Range("B9").Value = 'Hilary Osler'
$ws.Range("C7").Copy() | Out-Null
$ws.Range("C9").PasteSpecial(-4122) | Out-Null
$ws.Range("C9").Value = 'Innovation'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D9").PasteSpecial(-4122) | Out-Null
$ws.Range("D9").Value = 'Data Trends Tool'
$ws.Range("E7").Copy() | Out-Null
$ws.Range("E9").PasteSpecial(-4122) | Out-Null
$ws.Range("E9").Value = 'Yes'
$ws.Range("F7").Copy() | Out-Null
$ws.Range("F9").PasteSpecial(-4122) | Out-Null
$ws.Range("F9").Value = 'Mobile app'
$ws.Range("G7").Copy() | Out-Null
$ws.Range("G9").PasteSpecial(-4122) | Out-Null
$ws.Range("G9").Value = 'Maybe'
$ws.Range("H7").Copy() | Out-Null
$ws.Range("H9").PasteSpecial(-4122) | Out-Null
$ws.Range("H9").Value = 'The goal of this project is to create a model that uses previous climate data to predict future climate trends. A regression model will be used to generate predictions from past data.'
$ws.Range("I7").Copy() | Out-Null
$ws.Range("I9").PasteSpecial(-4122) | Out-Null
$ws.Range("I9").Value = 'A user-friendly app or website that demonstrates future trends (in the form of graphs/plots)  of climate data. '
$ws.Range("K7").Copy() | Out-Null
$ws.Range("K9").PasteSpecial(-4122) | Out-Null
$ws.Hyperlinks.Add($ws.Range("K9"), 'https://drive.google.com/open?id=1y1c5o-DKaeTKvwUwksnMdxWES1iHsB2R') | Out-Null
$ws.Range("K7").Copy() | Out-Null
$ws.Range("K9").PasteSpecial(-4122) | Out-Null

# Row 10
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A10").PasteSpecial(-4122) | Out-Null
$ws.Range("A10").Value = 44491.42447082176
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B10").PasteSpecial(-4122) | Out-Null
$ws.Range("B10").Value = 'Spencer Hill'
$ws.Range("C7").Copy() | Out-Null
$ws.Range("C10").PasteSpecial(-4122) | Out-Null
$ws.Range("C10").Value = 'Disruptive Tech'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4122) | Out-Null
$ws.Range("D10").Value = 'Novel Image generation using Quantum Generative Adversarial Networks.'
$ws.Range("E7").Copy() | Out-Null
$ws.Range("E10").PasteSpecial(-4122) | Out-Null
$ws.Range("E10").Value = 'No'
$ws.Range("H7").Copy() | Out-Null
$ws.Range("H10").PasteSpecial(-4122) | Out-Null
$ws.Range("H10").Value = 'We will be investigating the use of Quantum circuits to improve the performance of GAN Networks. Specifically, we will attempt to generate MNIST images using Quantum Computing. This project is more exploratory/research-based and I don''t think a product is within a reasonable scope of the project. It is already an ambitious project and after discussing with my general members they are more interested in further developing the architecture than integrating the code with a product (if we have time left over at the end). '

# Row 11
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A11").PasteSpecial(-4122) | Out-Null
$ws.Range("A11").Value = 44491.513493067134
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B11").PasteSpecial(-4122) | Out-Null
$ws.Range("B11").Value = 'Noah Cabral'
$ws.Range("C7").Copy() | Out-Null
$ws.Range("C11").PasteSpecial(-4122) | Out-Null
$ws.Range("C11").Value = 'DAIR'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D11").PasteSpecial(-4122) | Out-Null
$ws.Range("D11").Value = 'Hand Shape and Pose Tracking'
$ws.Range("E7").Copy() | Out-Null
$ws.Range("E11").PasteSpecial(-4122) | Out-Null
$ws.Range("E11").Value = 'No'
$ws.Range("H7").Copy() | Out-Null
$ws.Range("H11").PasteSpecial(-4122) | Out-Null
$ws.Range("H11").Value = 'The goal of the project is to build a machine learning model that can track the 3D mesh and pose of two hands. It should be invariant to occlusions from external objects as well as occlusions produced by inter- and intra-hand interactions. The motivation of this research is primarily an engineering one, where the goal is to develop an intelligent system that solves a real problem better than all alternative approaches. '

# Row 12
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A12").PasteSpecial(-4122) | Out-Null
$ws.Range("A12").Value = 44491.75316740741
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B12").PasteSpecial(-4122) | Out-Null
$ws.Range("B12").Value = 'Nicholas Murray'
$ws.Range("C7").Copy() | Out-Null
$ws.Range("C12").PasteSpecial(-4122) | Out-Null
$ws.Range("C12").Value = 'Innovation'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4122) | Out-Null
$ws.Range("D12").Value = 'Sports Betting App'
$ws.Range("E7").Copy() | Out-Null
$ws.Range("E12").PasteSpecial(-4122) | Out-Null
$ws.Range("E12").Value = 'Yes'
$ws.Range("F7").Copy() | Out-Null
$ws.Range("F12").PasteSpecial(-4122) | Out-Null
$ws.Range("F12").Value = 'Mobile app'
$ws.Range("G7").Copy() | Out-Null
$ws.Range("G12").PasteSpecial(-4122) | Out-Null
$ws.Range("G12").Value = 'Yes'
$ws.Range("H7").Copy() | Out-Null
$ws.Range("H12").PasteSpecial(-4122) | Out-Null
$ws.Range("H12").Value = 'The goal of the project is to develop an application which can be used by consumers to gain insights on what bets to place on the outcomes of sporting games. Like logistic regression. '
$ws.Range("I7").Copy() | Out-Null
$ws.Range("I12").PasteSpecial(-4122) | Out-Null
$ws.Range("I12").Value = 'A mobile application which can be used by consumers to gain insights on what bets to place on the outcomes of sporting games'
$ws.Range("K7").Copy() | Out-Null
$ws.Range("K12").PasteSpecial(-4122) | Out-Null
$ws.Hyperlinks.Add($ws.Range("K12"), 'https://drive.google.com/open?id=1G8WKWkEHy4i0yTjw_vzN3IUxesRelATO') | Out-Null
$ws.Range("K7").Copy() | Out-Null
$ws.Range("K12").PasteSpecial(-4122) | Out-Null

# Row 13
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null
$ws.Range("A13").Value = 44491.76838728009
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Range("B13").Value = 'Matt Wright'
$ws.Range("C7").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null
$ws.Range("C13").Value = 'Disruptive Tech'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4122) | Out-Null
$ws.Range("D13").Value = 'Variational Quantum Classifier'
$ws.Range("E7").Copy() | Out-Null
$ws.Range("E13").PasteSpecial(-4122) | Out-Null
$ws.Range("E13").Value = 'No'
$ws.Range("H7").Copy() | Out-Null
$ws.Range("H13").PasteSpecial(-4122) | Out-Null
$ws.Range("H13").Value = 'Image classification with a quantum computer'

# Row 14
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A14").PasteSpecial(-4122) | Out-Null
$ws.Range("A14").Value = 44491.81994216435
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B14").PasteSpecial(-4122) | Out-Null
$ws.Range("B14").Value = 'Braulio Antonio'
$ws.Range("C7").Copy() | Out-Null
$ws.Range("C14").PasteSpecial(-4122) | Out-Null
$ws.Range("C14").Value = 'Consulting'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4122) | Out-Null
$ws.Range("D14").Value = 'Evaluation of 3 industrial AI solutions'
$ws.Range("E7").Copy() | Out-Null
$ws.Range("E14").PasteSpecial(-4122) | Out-Null
$ws.Range("E14").Value = 'No'
$ws.Range("F7").Copy() | Out-Null
$ws.Range("F14").PasteSpecial(-4122) | Out-Null
$ws.Range("F14").Value = 'NA'
$ws.Range("G7").Copy() | Out-Null
$ws.Range("G14").PasteSpecial(-4122) | Out-Null
$ws.Range("G14").Value = 'No'
$ws.Range("H7").Copy() | Out-Null
$ws.Range("H14").PasteSpecial(-4122) | Out-Null
$ws.Range("H14").Value = 'NA'
$ws.Range("I7").Copy() | Out-Null
$ws.Range("I14").PasteSpecial(-4122) | Out-Null
$ws.Range("I14").Value = 'NA'
$ws.Range("J7").Copy() | Out-Null
$ws.Range("J14").PasteSpecial(-4122) | Out-Null
$ws.Range("J14").Value = 'NA'

# Row 15
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A15").PasteSpecial(-4122) | Out-Null
$ws.Range("A15").Value = 44491.88891119213
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4122) | Out-Null
$ws.Range("B15").Value = 'Camila Izquierdo'
$ws.Range("C7").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Range("C15").Value = 'Consulting'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null
$ws.Range("D15").Value = 'ADGA Research'
$ws.Range("E7").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null
$ws.Range("E15").Value = 'Yes'
$ws.Range("F7").Copy() | Out-Null
$ws.Range("F15").PasteSpecial(-4122) | Out-Null
$ws.Range("F15").Value = 'Website'
$ws.Range("G7").Copy() | Out-Null
$ws.Range("G15").PasteSpecial(-4122) | Out-Null
$ws.Range("G15").Value = 'No'
$ws.Range("H7").Copy() | Out-Null
$ws.Range("H15").PasteSpecial(-4122) | Out-Null
$ws.Range("H15").Value = 'Conducting a research study of different object detection and image classification architectures. We''ll be using a bunch of different models and just comparing them against each other using the client''s data to see which one gives the best result.'
$ws.Range("I7").Copy() | Out-Null
$ws.Range("I15").PasteSpecial(-4122) | Out-Null
$ws.Range("I15").Value = 'A website that could show the different architectures working alongside each other to classify an image.'
$ws.Range("K7").Copy() | Out-Null
$ws.Range("K15").PasteSpecial(-4122) | Out-Null
$ws.Hyperlinks.Add($ws.Range("K15"), 'https://drive.google.com/open?id=1QQFeoXX6NcoRiC_Y0016Ab-G0DVbPb-3') | Out-Null
$ws.Range("K7").Copy() | Out-Null
$ws.Range("K15").PasteSpecial(-4122) | Out-Null

# Row 16
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null
$ws.Range("A16").Value = 44491.88974216435
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B16").PasteSpecial(-4122) | Out-Null
$ws.Range("B16").Value = 'Ori Gurevich'
$ws.Range("C7").Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4122) | Out-Null
$ws.Range("C16").Value = 'Innovation'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4122) | Out-Null
$ws.Range("D16").Value = 'Legal Document Summarizer'
$ws.Range("E7").Copy() | Out-Null
$ws.Range("E16").PasteSpecial(-4122) | Out-Null
$ws.Range("E16").Value = 'Yes'
$ws.Range("F7").Copy() | Out-Null
$ws.Range("F16").PasteSpecial(-4122) | Out-Null
$ws.Range("F16").Value = 'Website'
$ws.Range("G7").Copy() | Out-Null
$ws.Range("G16").PasteSpecial(-4122) | Out-Null
$ws.Range("G16").Value = 'Maybe'
$ws.Range("H7").Copy() | Out-Null
$ws.Range("H16").PasteSpecial(-4122) | Out-Null
$ws.Range("H16").Value = 'Attempting to create a summarizer within certain fields '
$ws.Range("I7").Copy() | Out-Null
$ws.Range("I16").PasteSpecial(-4122) | Out-Null
$ws.Range("I16").Value = 'Website to upload legal documents, of which a summary is built and given back to user. Needs to adhere to privacy laws for obvious reasons'
$ws.Range("J7").Copy() | Out-Null
$ws.Range("J16").PasteSpecial(-4122) | Out-Null
$ws.Range("J16").Value = 'Unsure, likely frontend, maybe backend. Will remain in touch'
$ws.Range("K7").Copy() | Out-Null
$ws.Range("K16").PasteSpecial(-4122) | Out-Null
$ws.Hyperlinks.Add($ws.Range("K16"), 'https://drive.google.com/open?id=1zikUvWA5KO_DxoA_nqciaOxYIrks3Hst') | Out-Null
$ws.Range("K7").Copy() | Out-Null
$ws.Range("K16").PasteSpecial(-4122) | Out-Null

# Row 17
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A17").PasteSpecial(-4122) | Out-Null
$ws.Range("A17").Value = 44491.97619064814
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B17").PasteSpecial(-4122) | Out-Null
$ws.Range("B17").Value = 'Rabab Azeem'
$ws.Range("C7").Copy() | Out-Null
$ws.Range("C17").PasteSpecial(-4122) | Out-Null
$ws.Range("C17").Value = 'Consulting'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D17").PasteSpecial(-4122) | Out-Null
$ws.Range("D17").Value = 'Validere'
$ws.Range("E7").Copy() | Out-Null
$ws.Range("E17").PasteSpecial(-4122) | Out-Null
$ws.Range("E17").Value = 'Yes'
$ws.Range("F7").Copy() | Out-Null
$ws.Range("F17").PasteSpecial(-4122) | Out-Null
$ws.Range("F17").Value = 'Proof of concept that can be scaled (webapp)'
$ws.Range("G7").Copy() | Out-Null
$ws.Range("G17").PasteSpecial(-4122) | Out-Null
$ws.Range("G17").Value = 'Maybe'
$ws.Range("H7").Copy() | Out-Null
$ws.Range("H17").PasteSpecial(-4122) | Out-Null
$ws.Range("H17").Value = 'The project’s goal is to analyze documents released by companies and the government discussing the recent activities and trends about ESGs and sustainability initiatives. The project will be using NLP models to conduct topic modelling and sentiment analysis of sustainability methods used over time.'
$ws.Range("I7").Copy() | Out-Null
$ws.Range("I17").PasteSpecial(-4122) | Out-Null
$ws.Range("I17").Value = 'A website that can keep track of ESG trends and use the models we develop to display the current sentiments and ESG trends being discussed'
$ws.Range("J7").Copy() | Out-Null
$ws.Range("J17").PasteSpecial(-4122) | Out-Null
$ws.Range("J17").Value = 'Will have to ask the client about this'

# Row 18
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A18").PasteSpecial(-4122) | Out-Null
$ws.Range("A18").Value = 44492.0474934838
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B18").PasteSpecial(-4122) | Out-Null
$ws.Range("B18").Value = 'Courtney Orcutt'
$ws.Range("C7").Copy() | Out-Null
$ws.Range("C18").PasteSpecial(-4122) | Out-Null
$ws.Range("C18").Value = 'Consulting'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4122) | Out-Null
$ws.Range("D18").Value = 'AI in Fertility'
$ws.Range("E7").Copy() | Out-Null
$ws.Range("E18").PasteSpecial(-4122) | Out-Null
$ws.Range("E18").Value = 'No'
$ws.Range("G7").Copy() | Out-Null
$ws.Range("G18").PasteSpecial(-4122) | Out-Null
$ws.Range("G18").Value = 'No'
$ws.Range("H7").Copy() | Out-Null
$ws.Range("H18").PasteSpecial(-4122) | Out-Null
$ws.Range("H18").Value = 'We don''t yet have access to the data or even a very clear direction from the company as they are currently changing their wants based on the experience of the team. So will get back to you on the goals of the project since we may have to switch some things up.'
$ws.Range("K7").Copy() | Out-Null
$ws.Range("K18").PasteSpecial(-4122) | Out-Null
$ws.Hyperlinks.Add($ws.Range("K18"), 'https://drive.google.com/open?id=1eZY_vqSOn9e1wgS8Q_2yEn-m9r2yHn7t') | Out-Null
$ws.Range("K7").Copy() | Out-Null
$ws.Range("K18").PasteSpecial(-4122) | Out-Null


Write-Host "Rows 8-18 added successfully"
